$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks numeric need a forced-text / quote-prefix
# write followed by a style reset back to "Normal" so Excel keeps the literal
# text (e.g. "528.26") instead of silently converting it to a real number.
$ws.Range("D2").Value = "60.632.28"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.901.58"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.31%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "2.910.83"
$ws.Range("E9").Value = "  -2.08%  "
$ws.Range("E10").Value = "  -3.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.21%  "
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").Value = "3.408.98"
$ws.Range("E13").Value = "  -2.01%  "
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "60.620.09"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.50%  "
$ws.Range("D17").Value = "2.904.76"
$ws.Range("E17").Value = "  -2.06%  "
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("D30").Value = "0.0₃0861"
$ws.Range("E30").Value = "  -6.83%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "151.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.37%  "
$ws.Range("E35").Value = "  -5.03%  "
$ws.Range("E36").Value = "  -6.65%  "
$ws.Range("E38").Value = "  -5.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.26%  "
$ws.Range("D42").Value = "2.294.41"
$ws.Range("E42").Value = "  -4.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.648"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0582"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.75%  "
$ws.Range("E46").Value = "  +0.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.33%  "
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "251.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.00%  "
